$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D2:D51) keeps being stored as text, matching the
# original inline-string cells; otherwise Excel would auto-parse single-dot
# values like "310.77" as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.004.71"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "1.847.48"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  +0.48%  "
$ws.Range("D5").Value = "310.77"
$ws.Range("E5").Value = "  +1.29%  "
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Value = "0.4672"
$ws.Range("E7").Value = "  +3.17%  "
$ws.Range("D8").Value = "0.3636"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").Value = "0.07182"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").Value = "0.9353"
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").Value = "19.61"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "0.07684"
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").Value = "1.820.47"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").Value = "5.299"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "6.411"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "88.30"
$ws.Range("E16").Value = "  +2.99%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("E18").Value = "  +0.97%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "27.016.61"
$ws.Range("E20").Value = "  +1.98%  "
$ws.Range("D21").Value = "14.41"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "5.038"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("E23").Value = "  +1.24%  "
$ws.Range("D24").Value = "1.936"
$ws.Range("E24").Value = "  -1.51%  "
$ws.Range("D25").Value = "152.82"
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "2.031"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("D28").Value = "114.09"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("D29").Value = "4.938"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "0.08856"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").Value = "3.184"
$ws.Range("E31").Value = "  +2.60%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").Value = "2.851"
$ws.Range("E32").Value = "  +0.03%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "1.181"
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").Value = "0.7480"
$ws.Range("E34").Value = "  +2.81%  "
$ws.Range("D35").Value = "4.478"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").Value = "2.988"
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("D38").Value = "0.01940"
$ws.Range("E38").Value = "  +0.29%  "
$ws.Range("D39").Value = "0.05165"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").Value = "0.5138"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").Value = "6.916"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").Value = "0.1514"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("D43").Value = "8.194"
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("D44").Value = "10.47"
$ws.Range("E44").Value = "  +4.32%  "
$ws.Range("D45").Value = "0.4719"
$ws.Range("E45").Value = "  +0.66%  "
$ws.Range("E46").Value = "  +0.35%  "
$ws.Range("D47").Value = "100.39"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "1.606"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("D49").Value = "0.06049"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("D50").Value = "64.29"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "36.16"
